# Fix the codeforiati:category-name / group-code / group-name columns
# (E, F, G) in the SectorGroup sheet. The source data had the group-code
# value duplicated into the category-name column while the real
# category-name value had been pushed out into the group-name column.
# The correct mapping for every row (including the header row) is a
# left-rotation of the three columns:
#   new E (category-name) = old G
#   new F (group-code)    = old E
#   new G (group-name)    = old F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2
    $gVal = $gCell.Value2

    $eCell.Value2 = $gVal
    $fCell.Value2 = $eVal
    $gCell.Value2 = $fVal
}
